# Weekly update: insert a new price record as row 265 (pushing the
# existing rows 265-295 down to 266-296) on the Papa / Vega Monumental
# Concepcion sheet, and update the sheet's used-range dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 265; this shifts rows 265:295
# down to 266:296 and grows the used range to A1:R296.
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row 265 with the new weekly record.
$ws.Range("A265").Value = 11
$ws.Range("B265").Value = "Vega Monumental Concepción"
$ws.Range("C265").Value = "Bíobío"
$ws.Range("D265").Value = 44776
$ws.Range("E265").Value = 8
$ws.Range("F265").Value = 100114001
$ws.Range("G265").Value = "Papa"
$ws.Range("H265").Value = "Asterix"
$ws.Range("I265").Value = "1a (guarda lavada)"
$ws.Range("J265").Value = 300
$ws.Range("K265").Value = 7500
$ws.Range("L265").Value = 8000
$ws.Range("M265").Value = 7750
$ws.Range("N265").Value = "$/malla 25 kilos"
$ws.Range("O265").Value = "Región de La Araucanía"
$ws.Range("P265").Value = 310
$ws.Range("Q265").Value = 25
$ws.Range("R265").Value = "Hortaliza"
